$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new data row (row 12)
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "USA"
$ws.Range("C12").Value = "New York"
$ws.Range("D12").Value = "Kings"

# Fix capitalization typo: "jeffersen" -> "Jeffersen"
$ws.Range("D10").Value = "Jeffersen"

# Update selection to match target state
$ws.Range("E6").Select()
